$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting rows 103:216 down to 104:217
$ws.Rows("103:103").Insert()

# Populate the newly inserted row 103 with the new data record
$ws.Range("A103").Value = 3
$ws.Range("B103").Value = "Femacal de La Calera"
$ws.Range("C103").Value = "Coquimbo"
$ws.Range("D103").Value = Get-Date -Year 2023 -Month 1 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Range("E103").Value = 5
$ws.Range("F103").Value = 100112052
$ws.Range("G103").Value = "Albahaca"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 105
$ws.Range("K103").Value = 4000
$ws.Range("L103").Value = 4500
$ws.Range("M103").Value = 4238
$ws.Range("N103").Value = '$/docena de matas'
$ws.Range("O103").Value = "Provincia de Quillota"
$ws.Range("P103").Value = 706
$ws.Range("Q103").Value = 6
$ws.Range("R103").Value = "Hortaliza"

# Match the D-column date number format used by the rest of the column
$ws.Range("D103").NumberFormat = $ws.Range("D104").NumberFormat
